$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : update "propriétaire" / montants for the AV validation fix ---
$ws.Range("A2").Value = "001/TEST DR/AV"
$ws.Range("B2").Value = "Direction régionale"
# C2 is a purely-numeric string that must stay TEXT (quote-prefix keeps it textual
# without touching the shared default style used by every other cell).
$ws.Range("C2").Value = "'113564"
$ws.Range("D2").Value = "lala morale"
$ws.Range("E2").Value = "oui"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = "--"
$ws.Range("H2").Value = 24000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 24000
$ws.Range("M2").Value = 48000

# --- Row 3 : same treatment ---
$ws.Range("A3").Value = "001/TEST DR/AV"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "BB125874"
$ws.Range("D3").Value = "YASSINE TYEST"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = "--"
$ws.Range("H3").Value = 36000
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = "--"
$ws.Range("K3").Value = 5400
$ws.Range("L3").Value = 36000
$ws.Range("M3").Value = 66600

# --- Remove the now-obsolete rows 4-7 (used range shrinks to A1:M3) ---
$ws.Rows("4:7").Delete()
